# InsideBet Data: Automatizado
# The upstream fixture list advanced by one match: the earliest listed
# fixture (row 27 - Wk23, Fri 2026-02-20, Estrela vs Tondela) has fallen
# off the "upcoming" window, every subsequent fixture/row shifts up by
# one, and a new fixture for Wk34 (2026-05-17, Porto vs Santa Clara) that
# used to be last is now gone too (the table just got one row shorter
# overall, matching the new dimension A1:L144).
#
# Deleting the entire row 27 reproduces that shift precisely: Excel moves
# rows 28:145 up into 27:144 and the sheet's used range/dimension shrinks
# from A1:L145 to A1:L144 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("27").Delete()
